$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 125000360
$ws.Range("I28").Value = 125000360
$ws.Range("K28").Value = 125000360
$ws.Range("M28").Value = -124999875

$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").Value = $null

$ws.Range("H132").Value = 869.2
$ws.Range("I132").Value = 863.75
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 2591.25
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -61.25
$ws.Range("N132").Value = -8060

$ws.Range("H137").Value = 1402.7407
$ws.Range("I137").Value = 1299.2858
$ws.Range("J137").Value = 1764.8334
$ws.Range("K137").Value = 3897.8574
$ws.Range("L137").Value = 5294.5002
$ws.Range("M137").Value = -1347.8574
$ws.Range("N137").Value = -10394.5002

$ws.Range("H138").Value = 2026.3182
$ws.Range("I138").Value = 1813.2858
$ws.Range("J138").Value = 6500
$ws.Range("K138").Value = 5439.857400000001
$ws.Range("L138").Value = 19500
$ws.Range("M138").Value = -299.8574000000008
$ws.Range("N138").Value = -29780

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6329.135
$ws.Range("I32").Value = 3327.081
$ws.Range("K32").Value = 3327.081
$ws.Range("M32").Value = -3040.081

$ws.Range("H39").Value = 6497.5
$ws.Range("I39").Value = 6497.5
$ws.Range("K39").Value = 6497.5
$ws.Range("M39").Value = -5977.5

$ws.Range("H74").Value = 589.62067
$ws.Range("J74").Value = 1997
$ws.Range("L74").Value = 1997
$ws.Range("N74").Value = -3745

$ws.Range("H77").Value = 589.62067
$ws.Range("J77").Value = 1997
$ws.Range("L77").Value = 9985
$ws.Range("N77").Value = -18721

$ws.Range("H122").Value = 14426.538
$ws.Range("I122").Value = 7193.1816
$ws.Range("K122").Value = 21579.5448
$ws.Range("M122").Value = -19129.5448

$ws.Range("H132").Value = 3098.8333
$ws.Range("I132").Value = 3098.8333
$ws.Range("K132").Value = 9296.499899999999
$ws.Range("M132").Value = -6766.499899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 2995
$ws.Range("I33").Value = 2995
$ws.Range("K33").Value = 2995
$ws.Range("M33").Value = -2659

$ws.Range("H86").Value = 2798.6365
$ws.Range("J86").Value = 1500
$ws.Range("L86").Value = 1500
$ws.Range("N86").Value = -3746

$ws.Range("H89").Value = 2798.6365
$ws.Range("J89").Value = 1500
$ws.Range("L89").Value = 7500
$ws.Range("N89").Value = -18732

$ws.Range("H94").Value = 1067.8529
$ws.Range("I94").Value = 620.2222
$ws.Range("J94").Value = 2794.4285
$ws.Range("K94").Value = 620.2222
$ws.Range("L94").Value = 2794.4285
$ws.Range("M94").Value = -169.2222
$ws.Range("N94").Value = -3696.4285

$ws.Range("H107").Value = 1353.875
$ws.Range("I107").Value = 1187.6
$ws.Range("J107").Value = 1631
$ws.Range("K107").Value = 1187.6
$ws.Range("L107").Value = 1631
$ws.Range("M107").Value = 732.4000000000001
$ws.Range("N107").Value = -5471

$ws.Range("H134").Value = 2356.2222
$ws.Range("I134").Value = 1601
$ws.Range("K134").Value = 4803
$ws.Range("M134").Value = -2268

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5249.933
$ws.Range("J31").Value = 5943.1113
$ws.Range("L31").Value = 5943.1113
$ws.Range("N31").Value = -6533.1113

$ws.Range("H34").Value = 5249.933
$ws.Range("J34").Value = 5943.1113
$ws.Range("L34").Value = 5943.1113
$ws.Range("N34").Value = -6347.1113

$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").Value = $null

$ws.Range("H58").Value = 2350.5
$ws.Range("J58").Value = 4823.5713
$ws.Range("L58").Value = 4823.5713
$ws.Range("N58").Value = -5229.5713

$ws.Range("H132").Value = 1314.4
$ws.Range("I132").Value = 1314.4
$ws.Range("K132").Value = 3943.2
$ws.Range("M132").Value = -1413.2

$ws.Range("H136").Value = 2350.5
$ws.Range("J136").Value = 4823.5713
$ws.Range("L136").Value = 14470.7139
$ws.Range("N136").Value = -19570.7139

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 16666900
$ws.Range("I7").Value = 20000180
$ws.Range("K7").Value = 60000540
$ws.Range("M7").Value = -60000428

$ws.Range("H9").Value = 14000000
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").Value = $null

$ws.Range("H97").Value = 4792
$ws.Range("J97").Value = 4989.5
$ws.Range("L97").Value = 14968.5
$ws.Range("N97").Value = -15960.5

$ws.Range("H112").Value = 3150
$ws.Range("I112").Value = 2250
$ws.Range("K112").Value = 6750
$ws.Range("M112").Value = -5642

$ws.Range("H134").Value = 2489.75
$ws.Range("I134").Value = 2320
$ws.Range("J134").Value = 2999
$ws.Range("K134").Value = 6960
$ws.Range("L134").Value = 8997
$ws.Range("M134").Value = -1890
$ws.Range("N134").Value = -19137

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10000
$ws.Range("I70").Value = 10000
$ws.Range("K70").Value = 10000
$ws.Range("M70").Value = -9730

$ws.Range("H73").Value = 10000
$ws.Range("I73").Value = 10000
$ws.Range("K73").Value = 10000
$ws.Range("M73").Value = -9064

$ws.Range("H113").Value = 41709496
$ws.Range("I113").Value = 83359000
$ws.Range("K113").Value = 83359000
$ws.Range("M113").Value = -83356830

$ws.Range("H122").Value = 36476.758
$ws.Range("I122").Value = 1713.28
$ws.Range("J122").Value = 253748.5
$ws.Range("K122").Value = 5139.84
$ws.Range("L122").Value = 761245.5
$ws.Range("M122").Value = -2689.84
$ws.Range("N122").Value = -766145.5

$ws.Range("H132").Value = 1730.05
$ws.Range("I132").Value = 1481.1333
$ws.Range("J132").Value = 2476.8
$ws.Range("K132").Value = 4443.3999
$ws.Range("L132").Value = 7430.400000000001
$ws.Range("M132").Value = -1913.3999
$ws.Range("N132").Value = -12490.4

$ws.Range("H134").Value = 68493.5
$ws.Range("J134").Value = 68493.5
$ws.Range("L134").Value = 205480.5
$ws.Range("N134").Value = -210550.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 16250
$ws.Range("J64").Value = 16250
$ws.Range("L64").Value = 16250
$ws.Range("N64").Value = -16700

$ws.Range("H67").Value = 16250
$ws.Range("J67").Value = 16250
$ws.Range("L67").Value = 16250
$ws.Range("N67").Value = -17810

$ws.Range("H82").Value = 1503.7273
$ws.Range("I82").Value = 924
$ws.Range("K82").Value = 924
$ws.Range("M82").Value = -563

$ws.Range("H85").Value = 1503.7273
$ws.Range("I85").Value = 924
$ws.Range("K85").Value = 924
$ws.Range("M85").Value = 324

$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").Value = $null

$ws.Range("H132").Value = 4500
$ws.Range("I132").Value = 4500
$ws.Range("K132").Value = 13500
$ws.Range("M132").Value = -10970

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 249249.25
$ws.Range("I49").Value = 249333
$ws.Range("J49").Value = 248998
$ws.Range("K49").Value = 249333
$ws.Range("L49").Value = 248998
$ws.Range("M49").Value = -249103
$ws.Range("N49").Value = -249458

$ws.Range("H62").Value = 5443.5557
$ws.Range("J62").Value = 4750
$ws.Range("L62").Value = 4750
$ws.Range("N62").Value = -5998

$ws.Range("H63").Value = 14749.667
$ws.Range("J63").Value = 14749.667
$ws.Range("L63").Value = 14749.667
$ws.Range("N63").Value = -15997.667

$ws.Range("H65").Value = 5443.5557
$ws.Range("J65").Value = 4750
$ws.Range("L65").Value = 23750
$ws.Range("N65").Value = -29990

$ws.Range("H66").Value = 14749.667
$ws.Range("J66").Value = 14749.667
$ws.Range("L66").Value = 44249.001
$ws.Range("N66").Value = -50489.001

$ws.Range("H100").Value = 2565.05
$ws.Range("I100").Value = 2312
$ws.Range("J100").Value = 3155.5
$ws.Range("K100").Value = 4624
$ws.Range("L100").Value = 6311
$ws.Range("M100").Value = -4083
$ws.Range("N100").Value = -7393

$ws.Range("H132").Value = 1688.6471
$ws.Range("I132").Value = 1728.4375
$ws.Range("J132").Value = 1052
$ws.Range("K132").Value = 5185.3125
$ws.Range("L132").Value = 3156
$ws.Range("M132").Value = -2655.3125
$ws.Range("N132").Value = -8216

$ws.Range("H136").Value = 2183.3
$ws.Range("I136").Value = 972.6667
$ws.Range("K136").Value = 2918.0001
$ws.Range("M136").Value = -368.0001000000002
